$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 82) after the existing last row (row 81),
# matching the same layout used by every other data row:
#   A = date text, B = weekday text, C = hour (number), D = ranking (number)

# Column A holds a date-like string ("2025/10/09") that must stay literal
# text (as in every other row) rather than being auto-converted to a date
# serial number. A leading apostrophe forces Excel to store it as text;
# re-applying the "Normal" style afterwards strips the quote-prefix
# formatting so the cell ends up with the same (default) style as the
# rest of the data rows.
$ws.Range("A82").Value = "'2025/10/09"
$ws.Range("A82").Style = "Normal"

$ws.Range("B82").Value = "木"
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 201
